$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking price strings are not
# auto-converted to Excel numbers (the source data stores these as plain text).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.328.53"
$ws.Range("D3").Value = "1.935.46"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "0.7515"
$ws.Range("E5").Value = "  +5.46%  "
$ws.Range("D6").Value = "243.34"
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "28.10"
$ws.Range("E8").Value = "  +2.42%  "
$ws.Range("D9").Value = "0.3189"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "0.07270"
$ws.Range("E10").Value = "  +2.97%  "
$ws.Range("D11").Value = "0.7821"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").Value = "0.08044"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "1.968.75"
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "5.402"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "93.15"
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("D16").Value = "14.53"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "30.337.69"
$ws.Range("D18").Value = "6.108"
$ws.Range("E18").Value = "  +6.15%  "
$ws.Range("D19").Value = "252.33"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").Value = "0.000008064"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "2.187.63"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "6.703"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").Value = "9.593"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "165.07"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "19.11"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "0.1306"
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("D29").Value = "2.206"
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("D30").Value = "1.374"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").Value = "4.155"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").Value = "0.05301"
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("D35").Value = "1.335"
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("D36").Value = "0.7587"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("D37").Value = "2.791"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").Value = "0.01963"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "2.800"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "79.18"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").Value = "6.505"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("D42").Value = "0.4529"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").Value = "1.987"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "0.8394"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "7.740"
$ws.Range("E46").Value = "  +4.39%  "
$ws.Range("D47").Value = "10.06"
$ws.Range("E47").Value = "  +3.27%  "
$ws.Range("D48").Value = "101.74"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").Value = "37.67"
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("D50").Value = "0.1240"
$ws.Range("E50").Value = "  +9.15%  "
$ws.Range("D51").Value = "968.15"
$ws.Range("E51").Value = "  +5.47%  "
